$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the numeric "6" marker values in column E with the text value "n"
# for all rows where this applies.
$rows = @(4, 8, 10, 14, 20, 22, 25, 31, 35, 39, 41, 45)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 5).Value = "n"
}

# Update the selection to match the saved view state: whole column E selected,
# with E1 as the active cell.
$ws.Range("E1:E1048576").Select()
